$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New schedule data (rows 2-6), columns A-J:
# A=trialTrain, B=x_fixStart, C=y_fixStart, D=x_corrSteps, E=y_corrSteps,
# F=x_nrSteps, G=y_nrSteps, H=alienID, I=praclen, J=version
$data = @(
    @(1, 4, 7, 7, 4, 3, -3, 34, 5, "train_dim2_1"),
    @(2, 0, 5, 1, 0, 1, -5, 56, 5, "train_dim2_1"),
    @(3, 3, 6, 8, 5, 5, -1, 12, 5, "train_dim2_1"),
    @(4, 0, 7, 2, 3, 2, -4, 45, 5, "train_dim2_1"),
    @(5, 4, 6, 8, 4, 4, -2, 23, 5, "train_dim2_1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}
